# Greedy day heuristic update
# - Time_Window sheet: shift the two "busy" minute markers later in the day
#   (30 -> 570, 270 -> 810) for every row.
# - Add a new "Nurse_Type" sheet classifying each nurse index as RN or LVN.
# - Misc cursor/selection bookkeeping left behind by the edit session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Time_Window: move the non-zero block later in the day.
#    Every occurrence of 30 becomes 570, every occurrence of 270 becomes 810.
# ---------------------------------------------------------------------
$wsTW = $wb.Worksheets.Item("Time_Window")
for ($r = 2; $r -le 51; $r++) {
    for ($c = 2; $c -le 11; $c++) {
        $cell = $wsTW.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -eq 30) {
            $cell.Value = 570
        } elseif ($v -eq 270) {
            $cell.Value = 810
        }
    }
}

# ---------------------------------------------------------------------
# 2. Add the new Nurse_Type sheet (nurse index -> RN/LVN) after Min_Nurse.
# ---------------------------------------------------------------------
$wsMinNurse = $wb.Worksheets.Item("Min_Nurse")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNurseType = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNurseType.Name = "Nurse_Type"

# Header row
$wsNurseType.Range("A1").Value = "Index"
$wsNurseType.Range("B1").Value = "Type"

# Copy the "Index" column number formatting (bold + border + center/top)
# from Min_Nurse column A, which already uses the desired style.
$wsMinNurse.Range("A1:A51").Copy()
$wsNurseType.Range("A1:A51").PasteSpecial(-4122)

$rnIndexes = @(1,2,3,8,13,15,17,18,21,22,25,28,30,33,35,37,38,44,47,49)
$lvnIndexes = @(56,58,60,61,62,63,65,66,68,69,70,71,73,75,76,77,79,80,81,82,86,90,91,92,93,94,96,97,99,100)

$row = 2
foreach ($idx in $rnIndexes) {
    $wsNurseType.Cells.Item($row, 1).Value = $idx
    $wsNurseType.Cells.Item($row, 2).Value = "RN"
    $row = $row + 1
}
foreach ($idx in $lvnIndexes) {
    $wsNurseType.Cells.Item($row, 1).Value = $idx
    $wsNurseType.Cells.Item($row, 2).Value = "LVN"
    $row = $row + 1
}

$wsNurseType.Range("J28").Select()

# ---------------------------------------------------------------------
# 3. Leftover cursor/selection positions from the editing session.
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("H20").Select()

$wsHome = $wb.Worksheets.Item("C_home")
$wsHome.Columns("A").Select()

# Time_Window stays the active tab, cursor parked on M15.
$wsTW.Activate()
$wsTW.Range("M15").Select()

Write-Host "edit applied"
